# "Generate Report for Handback"
#
# This localization-status workbook is refreshed after a handback event:
#   - The Overview sheet's status columns flip from "Ready for handoff" to
#     "Handed back: in sync with en-US".
#   - The per-locale sheets (zh-cn, de-de) get their "Latest Target File" /
#     "Latest Handback File" / "Latest Handback DateTime" columns populated
#     for both data rows, with the new Target-File cell also turned into a
#     hyperlink (matching the existing hyperlink already used for the
#     Source File Name column).
#   - Column widths for the columns that now hold long file names / status
#     text are widened to fit.

$wb = $excel.ActiveWorkbook

$hyperlinkColor = 15570276  # RGB(100,149,237) == FF6495ED, the workbook's existing HyperLink font color

# ---------------------------------------------------------------------
# Overview sheet: status text for both rows/locale columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# zh-cn sheet: handback info for row 2 (0b82263a...) and row 3 (13fced72...)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1d49b6312a1b4711bc12962244b4275a648aab4/e2e/0b82263a-a77b-442e-9c12-54660209def2.md", [Type]::Missing, [Type]::Missing, "0b82263a-a77b-442e-9c12-54660209def2.md")
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("I2").Font.Color = $hyperlinkColor
$wsZh.Range("J2").Value = "0b82263a-a77b-442e-9c12-54660209def2.dc67295f6f2487c9bc2542d21970745352e0346f.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-31 02:08:40"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1d49b6312a1b4711bc12962244b4275a648aab4/e2e/13fced72-2230-4f4e-92d5-070dd3e9ba67.md", [Type]::Missing, [Type]::Missing, "13fced72-2230-4f4e-92d5-070dd3e9ba67.md")
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("I3").Font.Color = $hyperlinkColor
$wsZh.Range("J3").Value = "13fced72-2230-4f4e-92d5-070dd3e9ba67.a64300c2831f0e5f849c04dd8c9b50d6e9d79a7a.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-31 02:08:40"

$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(9).ColumnWidth = 39.15
$wsZh.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet: handback info for row 2 (0b82263a...) and row 3 (13fced72...)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1d49b6312a1b4711bc12962244b4275a648aab4/e2e/0b82263a-a77b-442e-9c12-54660209def2.md", [Type]::Missing, [Type]::Missing, "0b82263a-a77b-442e-9c12-54660209def2.md")
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("I2").Font.Color = $hyperlinkColor
$wsDe.Range("J2").Value = "0b82263a-a77b-442e-9c12-54660209def2.dc67295f6f2487c9bc2542d21970745352e0346f.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-31 02:08:57"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1d49b6312a1b4711bc12962244b4275a648aab4/e2e/13fced72-2230-4f4e-92d5-070dd3e9ba67.md", [Type]::Missing, [Type]::Missing, "13fced72-2230-4f4e-92d5-070dd3e9ba67.md")
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("I3").Font.Color = $hyperlinkColor
$wsDe.Range("J3").Value = "13fced72-2230-4f4e-92d5-070dd3e9ba67.a64300c2831f0e5f849c04dd8c9b50d6e9d79a7a.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-31 02:08:57"

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(9).ColumnWidth = 39.15
$wsDe.Columns.Item(10).ColumnWidth = 39.15
